# React-Django : Login And Notifications
# Append new TranslationsId rows (22=Email already registered, 23=Email/password
# incorrect, 24=Registration successful) as English/Turkish pairs, following the
# existing TextContentId / TranslationsId / LanguageId / Translations layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$rows = @(
    @(43, 22, 1, "Email is already registered."),
    @(44, 22, 2, "E-posta zaten kayıtlı."),
    @(45, 23, 1, "Email or password is incorrect."),
    @(46, 23, 2, "E-posta veya şifre hatalı."),
    @(47, 24, 1, "Registration Successful."),
    @(48, 24, 2, "Kayıt başarılı.")
)

$startRow = 44
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]
    $ws.Cells.Item($r, 1).Value = $data[0]
    $ws.Cells.Item($r, 2).Value = $data[1]
    $ws.Cells.Item($r, 3).Value = $data[2]
    $ws.Cells.Item($r, 4).Value = $data[3]
}

# Scroll/select to mirror the saved view state (top-left A24, active cell D48).
$excel.ActiveWindow.ScrollRow = 24
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D48").Select() | Out-Null
